$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$xlPasteFormats = -4122

# Template ranges used to copy the banded-table cell formatting used
# throughout the sheet.
$evenTemplate = $ws.Range("A14:C14")   # "even" banding style group (style ids 4,5,6)
$oddTemplate  = $ws.Range("A15:C15")   # "odd" banding style group  (style ids 7,8,9)

# Row 16 currently carries the special "last row" styling (style ids
# 10,11,12). Preserve a snapshot of that formatting in a scratch range
# before we touch row 16, because Range objects are live references and
# would otherwise reflect row 16's *new* styling once it's changed below.
$scratch = $ws.Range("A100:C100")
$ws.Range("A16:C16").Copy()
$scratch.PasteSpecial($xlPasteFormats)

# Row 16 is no longer the last row once new rows are appended, so restyle
# it as a normal (even) banded row.
$evenTemplate.Copy()
$ws.Range("A16:C16").PasteSpecial($xlPasteFormats)

# New survey responses to append to the table.
$data = @(
    @(45755.591476504633, 350, 870),
    @(45755.609441886569, 800, 500),
    @(45755.657248310185, 450, 120),
    @(45755.657509606477, 600, 1500),
    @(45755.779672268516, 500, 600),
    @(45755.779805428239, 100, 50),
    @(45755.779920405097, 350, 350),
    @(45755.780091770837, 500, 680)
)

$count = $data.Count
for ($i = 0; $i -lt $count; $i++) {
    $d = $data[$i]
    $newRow = $lo.ListRows.Add()
    $r = $newRow.Range
    $r.Cells.Item(1,1).Value = $d[0]
    $r.Cells.Item(1,2).Value = $d[1]
    $r.Cells.Item(1,3).Value = $d[2]

    $isLast = ($i -eq ($count - 1))
    if ($isLast) {
        # Final row gets the preserved "last row" styling.
        $scratch.Copy()
    } elseif (($i % 2) -eq 0) {
        # i=0 -> row17 (odd row number) -> odd banding
        $oddTemplate.Copy()
    } else {
        # i=1 -> row18 (even row number) -> even banding
        $evenTemplate.Copy()
    }
    $r.PasteSpecial($xlPasteFormats)
}

# Remove the scratch helper range so it doesn't linger in the sheet.
$scratch.Clear()

# Update the active selection to reflect the new data entry point, as in
# the target workbook.
[void]$ws.Range("B2").Select()

Write-Host "Final table range:" $lo.Range.Address()
